# Add the new 조의금(condolence money) entries to the bottom of Sheet1's table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Name (column A) / amount in 만원 units (column B) pairs to append,
# continuing directly after the existing last row (307).
$newEntries = @(
    @("김장식", 5),
    @("이학진", 5),
    @("조현준", 5),
    @("이상조", 5),
    @("안승호", 5),
    @("양덕우", 10),
    @("박진욱", 10),
    @("조합골프 총무", 10),
    @("채일교", 10)
)

$startRow = 308
for ($i = 0; $i -lt $newEntries.Count; $i++) {
    $row = $startRow + $i
    $name = $newEntries[$i][0]
    $amount = $newEntries[$i][1]
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $amount
}

# Move the active selection to the last newly-added cell, matching the
# author's final cursor position after pasting the new rows in.
$ws.Range("A316").Select()
